$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.351.49"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "2.605.63"
$ws.Range("E3").Value = "  +3.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.93"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.51"
$ws.Range("E6").Value = "  +4.63%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.12"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0818"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "3.004.68"
$ws.Range("E13").Value = "  +3.71%  "
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "2.644.14"
$ws.Range("E15").Value = "  +5.07%  "
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.852"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "43.440.16"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("E19").Value = "  +3.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.79"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.73"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.27"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.99"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("E25").Value = "  +4.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.35"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.49"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.35"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.90"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.69"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("E33").Value = "  +7.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.18"
$ws.Range("E34").Value = "  +3.79%  "
$ws.Range("E35").Value = "  +4.19%  "
$ws.Range("E36").Value = "  +3.06%  "
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("E39").Value = "  +9.74%  "
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.92"
$ws.Range("E41").Value = "  -3.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.03"
$ws.Range("E42").Value = "  +7.28%  "
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").Value = "2.014.57"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.02"
$ws.Range("E47").Value = "  +2.53%  "
$ws.Range("D48").Value = "2.854.33"
$ws.Range("E48").Value = "  +3.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "83.87"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.25"
$ws.Range("E50").Value = "  +2.99%  "
$ws.Range("E51").Value = "  +4.74%  "
